# Shift all timestamps in column A (rows 2-97) forward by 22 days,
# and update the solar production values in column B for rows 30-40
# to reflect the new day's generation curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2() + 22
}

$bUpdates = @{
    30 = 1
    31 = 18
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
}

foreach ($r in $bUpdates.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $bUpdates[$r]
}
